# edit.ps1 - applies the STRIDE "Spoofing" table cleanup described in the
# commit diff:
#   1. Split "ARP Spoofing, IP spoofing, DNS spoofing " into two runs and
#      drop "DNS spoofing" from the visible text.
#   2. Fill in the (previously empty) mitigation cell for "Spoofing a
#      machine" with two mitigation paragraphs.
#   3. Remove the whole "Spoofing a process" row (its content moved/merged
#      elsewhere).
#   4. Re-balance the <w:lastRenderedPageBreak/> markers that shift because
#      of the row deletions: one now lands on "Repudiating an action" and
#      another on "Inject a command", while the one that used to sit in
#      front of "Information Disclosure" (and the extra blank paragraph
#      before it) goes away.
#   5. Remove the whole "Spoofing a user" row.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphXml($range, [string]$innerXml) {
    # Replace the *entire* contents of the paragraph(s) covered by $range
    # (start..end-1, i.e. excluding the trailing paragraph mark) with a
    # brand-new <w:p> fragment. Operating on the whole paragraph (not just
    # a sub-run) avoids InsertXML silently dropping sibling runs.
    $target = $d.Range($range.Start, $range.End - 1)
    $target.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

# ---------------------------------------------------------------------
# 1 & 2: the "Spoofing a machine" row in the Spoofing table.
# ---------------------------------------------------------------------
$spoofTable = $d.Tables.Item(1)

# Col 2: "ARP Spoofing, IP spoofing, DNS spoofing " -> two runs, dropping
# ", DNS spoofing" and adding a closing ". ".
$machineAttack = $spoofTable.Cell(3, 2).Range
Set-ParagraphXml $machineAttack (
    '<w:r><w:t>ARP Spoofing, IP spoofing</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">. </w:t></w:r>'
)

# Col 4: was a single empty paragraph; becomes two mitigation paragraphs
# plus trailing blank paragraphs.
$machineMitigation = $spoofTable.Cell(3, 4).Range
$target = $d.Range($machineMitigation.Start, $machineMitigation.End - 1)
$target.InsertXML(
    "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Implement spoofing detection software. </w:t></w:r></w:p>" +
    "<w:p $wNs/>" +
    "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Implement protocols such as HTTPs, that will drastically reduce the chance of a successful ARP poisoning attack and mitigate IP spoofing. </w:t></w:r></w:p>" +
    "<w:p $wNs/>" +
    "<w:p $wNs/>"
)

# ---------------------------------------------------------------------
# 3: delete the "Spoofing a process" row entirely (row 4 of the same
# table - row indices don't shift from the edits above).
# ---------------------------------------------------------------------
$spoofTable.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# 4a: "Repudiating an action" now picks up a lastRenderedPageBreak.
# ---------------------------------------------------------------------
$repudiationTable = $d.Tables.Item(3)
$repudiateCell = $repudiationTable.Cell(2, 1).Range
Set-ParagraphXml $repudiateCell (
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Repudiating an action </w:t></w:r>'
)

# ---------------------------------------------------------------------
# 4b: drop the lastRenderedPageBreak from "Information Disclosure" and
# collapse the two blank paragraphs above it down to one.
# ---------------------------------------------------------------------
$headingSearch = $d.Content
$headingSearch.Find.Execute("Information Disclosure:") | Out-Null
$headingStart = $headingSearch.Start
$headingPara = $d.Range($headingStart, $headingStart + 1).Paragraphs.Item(1).Range

Set-ParagraphXml $headingPara (
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Information Disclosure</w:t></w:r>' +
    '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>:</w:t></w:r>'
)

# Remove one of the two blank paragraphs that precede the heading (the
# pair sits right before $headingStart).
$d.Range($headingStart - 2, $headingStart - 1).Delete() | Out-Null

# ---------------------------------------------------------------------
# 4c: "Inject a command" now picks up a lastRenderedPageBreak.
# ---------------------------------------------------------------------
$injectTable = $d.Tables.Item(6)
$injectCell = $injectTable.Cell(3, 1).Range
Set-ParagraphXml $injectCell (
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Inject a command  </w:t></w:r>'
)

# ---------------------------------------------------------------------
# 5: delete the "Spoofing a user" row entirely (row 4 of the same table).
# ---------------------------------------------------------------------
$injectTable.Rows.Item(4).Delete()

Write-Output "edit.ps1 completed"
